$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "217.60", "1.01", "7.70").
# Pre-format as Text so Excel keeps the exact literal string instead of silently
# parsing/rounding it into a float (which would also flip the stored cell type).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.852.23"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.646.89"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +1.16%  "
$ws.Range("D5").Value = "217.60"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "19.23"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "1.871.85"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "1.656.93"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "64.87"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").Value = "26.849.71"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "0.0$([char]0x2083)0739"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").Value = "214.16"
$ws.Range("E19").Value = "  -4.23%  "
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").Value = "2.41"
$ws.Range("E22").Value = "  +11.61%  "
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "9.37"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").Value = "145.33"
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("D28").Value = "7.09"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("D34").Value = "1.277.33"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("E37").Value = "  -3.62%  "
$ws.Range("D38").Value = "0.537"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").Value = "0.828"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "0.814"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("D44").Value = "1.797.78"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "91.70"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("D46").Value = "58.93"
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0$([char]0x2086)0103"
$ws.Range("E48").Value = "  -2.49%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0519"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.70"
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0976"
$ws.Range("E51").Value = "  -0.52%  "
